# Estacion.h.docx edit script
#
# Adds:
#   - a forward declaration of class Red before "class Estacion {"
#   - a "friend class Red;" declaration (with comment) before "private:"
#   - three new accessor method declarations (with comments) before the
#     closing "};" of class Estacion
#   - a trailing blank paragraph right after "#endif"
#
# Strategy: locate a stable anchor paragraph with Find, remember its
# paragraph Index (n) *before* inserting, insert the required number of
# blank paragraphs right before it (InsertParagraphBefore always inserts
# immediately above the anchor, pushing the anchor further down), then
# fill in the text of the freshly-created blank paragraphs using plain
# absolute indices n, n+1, n+2, ... via $d.Paragraphs.Item(...).

$d = $word.ActiveDocument

function Insert-ParagraphsBefore($anchorText, $newTexts) {
    $range = $d.Content
    $range.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $n = $range.Paragraphs.Item(1).Index
    $count = $newTexts.Count
    for ($k = 0; $k -lt $count; $k++) {
        $range.InsertParagraphBefore()
    }
    for ($k = 0; $k -lt $count; $k++) {
        $txt = $newTexts[$k]
        if ($txt -ne "") {
            $d.Paragraphs.Item($n + $k).Range.Text = $txt
        }
    }
}

# 1) Before "class Estacion {"
Insert-ParagraphsBefore "class Estacion {" @(
    "class Red; // Declaración anticipada de la clase Red",
    ""
)

# 2) Before "private:"
Insert-ParagraphsBefore "private:" @(
    "    // Declaración de amistad con la clase Red",
    "    friend class Red;",
    ""
)

# 3) Before the closing "};" of class Estacion
Insert-ParagraphsBefore "};" @(
    "",
    "    // Método para obtener el siguiente nodo",
    "    Estacion* getSiguienteEstacion() const;",
    "",
    "    // Método para obtener el nodo anterior",
    "    Estacion* getAnteriorEstacion() const;",
    "",
    "    // Método público para obtener el puntero siguienteEstacion",
    "    Estacion* obtenerSiguienteEstacion() const;"
)

# 4) A blank paragraph right after "#endif"
$range = $d.Content
$range.Find.Execute("#endif", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$range.InsertParagraphAfter()

Write-Output "done"
